$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2339.3076
$ws.Range("I98").Value = 1249.7273
$ws.Range("J98").Value = 8332
$ws.Range("K98").Value = 1249.7273
$ws.Range("L98").Value = 8332
$ws.Range("M98").Value = 248.2727
$ws.Range("N98").Value = -11328
$ws.Range("H113").Value = 1974.2727
$ws.Range("I113").Value = 1448.3334
$ws.Range("K113").Value = 1448.3334
$ws.Range("M113").Value = 1805.6666
$ws.Range("H122").Value = 2339.3076
$ws.Range("I122").Value = 1249.7273
$ws.Range("J122").Value = 8332
$ws.Range("K122").Value = 3749.1819
$ws.Range("L122").Value = 24996
$ws.Range("M122").Value = -1299.1819
$ws.Range("N122").Value = -29896
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H135").Value = 1813.8214
$ws.Range("I135").Value = 1137.8182
$ws.Range("K135").Value = 10240.3638
$ws.Range("M135").Value = -7705.363799999999

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2650.2856
$ws.Range("I2").Value = 2258.6667
$ws.Range("K2").Value = 2258.6667
$ws.Range("M2").Value = -2145.6667
$ws.Range("H61").Value = 2786.9
$ws.Range("I61").Value = 1763.2222
$ws.Range("K61").Value = 1763.2222
$ws.Range("M61").Value = -1551.2222
$ws.Range("H102").Value = 2938.8823
$ws.Range("I102").Value = 2497.5625
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 2497.5625
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = -875.5625
$ws.Range("N102").Value = -13244
$ws.Range("H110").Value = 2614.5557
$ws.Range("I110").Value = 2645.4285
$ws.Range("K110").Value = 2645.4285
$ws.Range("M110").Value = -600.4285
$ws.Range("H116").Value = 2650.2856
$ws.Range("I116").Value = 2258.6667
$ws.Range("K116").Value = 2258.6667
$ws.Range("M116").Value = 35.33329999999978
$ws.Range("H122").Value = 7425.9
$ws.Range("I122").Value = 4332.375
$ws.Range("K122").Value = 12997.125
$ws.Range("M122").Value = -10547.125
$ws.Range("H136").Value = 2786.9
$ws.Range("I136").Value = 1763.2222
$ws.Range("K136").Value = 5289.6666
$ws.Range("M136").Value = -2739.6666

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2650.2856
$ws.Range("I3").Value = 2258.6667
$ws.Range("K3").Value = 2258.6667
$ws.Range("M3").Value = -2144.6667
$ws.Range("H64").Value = 448.44446
$ws.Range("I64").Value = 324.5
$ws.Range("J64").Value = 483.85715
$ws.Range("K64").Value = 324.5
$ws.Range("L64").Value = 483.85715
$ws.Range("M64").Value = -99.5
$ws.Range("N64").Value = -933.85715
$ws.Range("H67").Value = 448.44446
$ws.Range("I67").Value = 324.5
$ws.Range("J67").Value = 483.85715
$ws.Range("K67").Value = 324.5
$ws.Range("L67").Value = 483.85715
$ws.Range("M67").Value = 455.5
$ws.Range("N67").Value = -2043.85715
$ws.Range("H80").Value = 13213.1875
$ws.Range("J80").Value = 779.4167
$ws.Range("L80").Value = 779.4167
$ws.Range("N80").Value = -2775.4167
$ws.Range("H83").Value = 13213.1875
$ws.Range("J83").Value = 779.4167
$ws.Range("L83").Value = 3897.0835
$ws.Range("N83").Value = -13881.0835

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 160607
$ws.Range("I74").Value = 24500
$ws.Range("J74").Value = 183291.5
$ws.Range("K74").Value = 24500
$ws.Range("L74").Value = 183291.5
$ws.Range("M74").Value = -23626
$ws.Range("N74").Value = -185039.5
$ws.Range("H77").Value = 160607
$ws.Range("I77").Value = 24500
$ws.Range("J77").Value = 183291.5
$ws.Range("K77").Value = 73500
$ws.Range("L77").Value = 549874.5
$ws.Range("M77").Value = -69132
$ws.Range("N77").Value = -558610.5
$ws.Range("H105").Value = 2835.3794
$ws.Range("J105").Value = 3577
$ws.Range("L105").Value = 3577
$ws.Range("N105").Value = -7071
$ws.Range("H112").Value = 171392.72
$ws.Range("J112").Value = 171392.72
$ws.Range("L112").Value = 171392.72
$ws.Range("N112").Value = -174346.72

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 163
$ws.Range("I13").Value = 151.33333
$ws.Range("K13").Value = 453.99999
$ws.Range("M13").Value = -285.99999
$ws.Range("H17").Value = 291.94736
$ws.Range("I17").Value = 200.6
$ws.Range("J17").Value = 393.44446
$ws.Range("K17").Value = 601.8
$ws.Range("L17").Value = 1180.33338
$ws.Range("M17").Value = -432.8
$ws.Range("N17").Value = -1518.33338
$ws.Range("H46").Value = 600.875
$ws.Range("I46").Value = 203.25
$ws.Range("J46").Value = 998.5
$ws.Range("K46").Value = 609.75
$ws.Range("L46").Value = 2995.5
$ws.Range("M46").Value = -518.75
$ws.Range("N46").Value = -3177.5
$ws.Range("H107").Value = 1291.2941
$ws.Range("J107").Value = 1281.3334
$ws.Range("L107").Value = 3844.0002
$ws.Range("N107").Value = -7684.0002
$ws.Range("H109").Value = 527.125
$ws.Range("I109").Value = 248.4
$ws.Range("J109").Value = 991.6667
$ws.Range("K109").Value = 745.2
$ws.Range("L109").Value = 2975.0001
$ws.Range("M109").Value = 294.8
$ws.Range("N109").Value = -5055.0001
$ws.Range("H113").Value = 1123.6364
$ws.Range("J113").Value = 1123.6364
$ws.Range("L113").Value = 3370.9092
$ws.Range("N113").Value = -7710.9092
$ws.Range("H119").Value = 6697
$ws.Range("I119").Value = 4663.3335
$ws.Range("K119").Value = 13990.0005
$ws.Range("M119").Value = -9152.000499999998
$ws.Range("H122").Value = 789
$ws.Range("I122").Value = 756
$ws.Range("J122").Value = 834.375
$ws.Range("K122").Value = 6804
$ws.Range("L122").Value = 7509.375
$ws.Range("M122").Value = -4354
$ws.Range("N122").Value = -12409.375

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2995.625
$ws.Range("I102").Value = 1900.3334
$ws.Range("K102").Value = 1900.3334
$ws.Range("M102").Value = -278.3334
$ws.Range("H122").Value = 2793.0605
$ws.Range("I122").Value = 1193.2916
$ws.Range("K122").Value = 3579.8748
$ws.Range("M122").Value = -1129.8748
$ws.Range("H126").Value = 4344.154
$ws.Range("I126").Value = 3330.4443
$ws.Range("K126").Value = 9991.332900000001
$ws.Range("M126").Value = -7521.332900000001

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5680.3335
$ws.Range("I61").Value = 4538.0435
$ws.Range("J61").Value = 12248.5
$ws.Range("K61").Value = 4538.0435
$ws.Range("L61").Value = 12248.5
$ws.Range("M61").Value = -4336.0435
$ws.Range("N61").Value = -12652.5
$ws.Range("H113").Value = 5680.3335
$ws.Range("I113").Value = 4538.0435
$ws.Range("J113").Value = 12248.5
$ws.Range("K113").Value = 4538.0435
$ws.Range("L113").Value = 12248.5
$ws.Range("M113").Value = -2368.0435
$ws.Range("N113").Value = -16588.5
$ws.Range("H136").Value = 7566.3193
$ws.Range("I136").Value = 2444.5715
$ws.Range("K136").Value = 7333.7145
$ws.Range("M136").Value = -4783.7145

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 1725
$ws.Range("I13").Value = 1725
$ws.Range("K13").Value = 1725
$ws.Range("M13").Value = -1585
$ws.Range("H75").Value = 152484.88
$ws.Range("I75").Value = 59999.5
$ws.Range("J75").Value = 183313.33
$ws.Range("K75").Value = 59999.5
$ws.Range("L75").Value = 183313.33
$ws.Range("M75").Value = -59063.5
$ws.Range("N75").Value = -185185.33
$ws.Range("H78").Value = 152484.88
$ws.Range("I78").Value = 59999.5
$ws.Range("J78").Value = 183313.33
$ws.Range("K78").Value = 179998.5
$ws.Range("L78").Value = 549939.99
$ws.Range("M78").Value = -175318.5
$ws.Range("N78").Value = -559299.99
$ws.Range("H113").Value = 418.13794
$ws.Range("I113").Value = 288.5263
$ws.Range("K113").Value = 865.5789
$ws.Range("M113").Value = 1304.4211

